$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.582.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.73%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.299.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.02%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.25%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'541.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.18%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'128.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -4.31%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.20%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.568"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -3.54%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'2.297.36"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.87%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.100"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.75%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'  -0.73%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'  -0.76%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "'  -2.14%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'23.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -4.42%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("B15").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'2.707.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -2.45%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("B16").Value = "'WrappedBTC"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'59.422.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.90%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = "'  -2.23%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'2.291.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.80%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'10.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -3.27%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'4.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -4.80%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'309.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.43%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'  -3.48%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.68%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'62.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.56%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = "'  -3.44%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'  +0.21%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'7.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -4.58%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'1.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.73%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "'  +2.21%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'171.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.04%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("E31").Value = "'  -2.18%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'0.0₃0711"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -5.55%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'5.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.82%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'0.377"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -3.19%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "'  -0.02%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  -8.34%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'17.57"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.30%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.04%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'3.97"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -3.76%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'312.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.27%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'37.45"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -2.20%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'1.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -3.56%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'135.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -5.97%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'3.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -2.12%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.0938"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -2.20%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Value = "'  +0.06%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.0₆0231"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +27.21%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'18.46"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.43%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = "'  -2.60%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'  -0.67%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Value = "'  -0.22%  "
$ws.Range("E51").Style = "Normal"

